$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# '77÷7=11, 0' -> '73÷2=36, 1'
$t.Cell(1, 1).Range.Text = "73÷2=36, 1"

# '34÷9=3, 7' -> '58÷4=14, 2'
$t.Cell(1, 2).Range.Text = "58÷4=14, 2"

# '38÷8=4, 6' -> '45÷5=9, 0'
$t.Cell(1, 3).Range.Text = "45÷5=9, 0"

# '74÷4=18, 2' -> '35÷7=5, 0'
$t.Cell(1, 4).Range.Text = "35÷7=5, 0"

# '74÷9=8, 2' -> '51÷8=6, 3'
$t.Cell(1, 5).Range.Text = "51÷8=6, 3"

# '72÷9=8, 0' -> '70÷6=11, 4'
$t.Cell(5, 1).Range.Text = "70÷6=11, 4"

# '61÷6=10, 1' -> '47÷4=11, 3'
$t.Cell(5, 2).Range.Text = "47÷4=11, 3"

# '56÷4=14, 0' -> '29÷9=3, 2'
$t.Cell(5, 3).Range.Text = "29÷9=3, 2"

# '71÷5=14, 1' -> '97÷6=16, 1'
$t.Cell(5, 4).Range.Text = "97÷6=16, 1"

# '35÷4=8, 3' -> '51÷8=6, 3'
$t.Cell(5, 5).Range.Text = "51÷8=6, 3"

# '52÷3=17, 1' -> '84÷8=10, 4'
$t.Cell(9, 1).Range.Text = "84÷8=10, 4"

# '60÷5=12, 0' -> '74÷6=12, 2'
$t.Cell(9, 2).Range.Text = "74÷6=12, 2"

# '64÷5=12, 4' -> '90÷4=22, 2'
$t.Cell(9, 3).Range.Text = "90÷4=22, 2"

# '44÷5=8, 4' -> '46÷3=15, 1'
$t.Cell(9, 4).Range.Text = "46÷3=15, 1"

# '89÷2=44, 1' -> '21÷7=3, 0'
$t.Cell(9, 5).Range.Text = "21÷7=3, 0"

# '14÷4=3, 2' -> '74÷7=10, 4'
$t.Cell(13, 1).Range.Text = "74÷7=10, 4"

# '97÷3=32, 1' -> '92÷3=30, 2'
$t.Cell(13, 2).Range.Text = "92÷3=30, 2"

# '17÷8=2, 1' -> '13÷6=2, 1'
$t.Cell(13, 3).Range.Text = "13÷6=2, 1"

# '33÷7=4, 5' -> '31÷5=6, 1'
$t.Cell(13, 4).Range.Text = "31÷5=6, 1"

# '88÷8=11, 0' -> '61÷7=8, 5'
$t.Cell(13, 5).Range.Text = "61÷7=8, 5"

# '59÷8=7, 3' -> '17÷8=2, 1'
$t.Cell(17, 1).Range.Text = "17÷8=2, 1"

# '48÷6=8, 0' -> '11÷9=1, 2'
$t.Cell(17, 2).Range.Text = "11÷9=1, 2"

# '96÷4=24, 0' -> '64÷8=8, 0'
$t.Cell(17, 3).Range.Text = "64÷8=8, 0"

# '20÷7=2, 6' -> '95÷7=13, 4'
$t.Cell(17, 4).Range.Text = "95÷7=13, 4"

# '65÷9=7, 2' -> '36÷2=18, 0'
$t.Cell(17, 5).Range.Text = "36÷2=18, 0"
